$d = $word.ActiveDocument

# --- Step 1: append the student's name to the "Nama :" line ---
$namaPara = $d.Paragraphs.Item(4)
$namaRange = $namaPara.Range
$namaRange.End = $namaRange.End - 1   # exclude the paragraph mark
$namaRange.Collapse(0)                # wdCollapseEnd
$namaRange.InsertAfter([string]::Concat("Daffa Cesario Safi", [char]0x2019, "i"))

# --- Step 2: insert a brand-new paragraph right after it ---
$namaPara2 = $d.Paragraphs.Item(4)
$afterNamaRange = $namaPara2.Range
$afterNamaRange.InsertParagraphAfter()

# --- Step 3: fill the new paragraph with the submission-link text ---
$linkPara = $d.Paragraphs.Item(5)
$linkRange = $linkPara.Range
$linkRange.End = $linkRange.End - 1
$linkRange.Text = "Link pengumpulan tugas : "

$linkPara2 = $d.Paragraphs.Item(5)
$linkRange2 = $linkPara2.Range
$linkRange2.End = $linkRange2.End - 1
$linkRange2.Collapse(0)
$linkRange2.InsertAfter("https://github.com/daffacesario46/tugasdatabase.git")
